$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 5.692099502487557
$ws.Cells.Item(2, 2).Value = 0.008
$ws.Cells.Item(2, 3).Value = 21.53199999999999

$ws.Cells.Item(3, 1).Value = 1.419999999999999
$ws.Cells.Item(3, 3).Value = 6.300000000000003

$ws.Cells.Item(4, 1).Value = 1.913353233830843
$ws.Cells.Item(4, 3).Value = 8.071999999999997

$ws.Cells.Item(5, 1).Value = 1.708378109452733
$ws.Cells.Item(5, 3).Value = 7.228000000000003

$ws.Cells.Item(6, 1).Value = 3.659502487562182
$ws.Cells.Item(6, 2).Value = 0.004
$ws.Cells.Item(6, 3).Value = 14.836

$ws.Cells.Item(7, 1).Value = 5.053592039800984
$ws.Cells.Item(7, 2).Value = 0.044
$ws.Cells.Item(7, 3).Value = 19.928

$ws.Cells.Item(8, 1).Value = 7.040199004975123
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 26.95999999999999

$ws.Cells.Item(9, 1).Value = 3.267840796019891
$ws.Cells.Item(9, 3).Value = 14.34

$ws.Cells.Item(10, 1).Value = 6.658646766169147
$ws.Cells.Item(10, 2).Value = 0.02
$ws.Cells.Item(10, 3).Value = 26.144

$ws.Cells.Item(11, 1).Value = 6.048039800995014
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 22.10800000000001

$ws.Cells.Item(12, 1).Value = 4.655641791044771
$ws.Cells.Item(12, 2).Value = 0.024
$ws.Cells.Item(12, 3).Value = 18.80399999999999

$ws.Cells.Item(13, 1).Value = 6.695383084577109
$ws.Cells.Item(13, 2).Value = 0.06000000000000002
$ws.Cells.Item(13, 3).Value = 25.96399999999999

$ws.Cells.Item(14, 1).Value = 6.718905472636805
$ws.Cells.Item(14, 2).Value = 0.06000000000000001
$ws.Cells.Item(14, 3).Value = 26.41599999999999

$ws.Cells.Item(15, 1).Value = 2.131860696517406
$ws.Cells.Item(15, 3).Value = 9.131999999999998

$ws.Cells.Item(16, 1).Value = 3.663860696517405
$ws.Cells.Item(16, 3).Value = 15.264

$ws.Cells.Item(17, 1).Value = 4.925711442786061
$ws.Cells.Item(17, 2).Value = 0.1040000000000001
$ws.Cells.Item(17, 3).Value = 20.084

$ws.Cells.Item(18, 1).Value = 1.990328358208951
$ws.Cells.Item(18, 3).Value = 8.735999999999997

$ws.Cells.Item(19, 1).Value = 6.196019900497503
$ws.Cells.Item(19, 2).Value = 0.016
$ws.Cells.Item(19, 3).Value = 26.008

$ws.Cells.Item(20, 1).Value = 3.870845771144271
$ws.Cells.Item(20, 2).Value = 0.044
$ws.Cells.Item(20, 3).Value = 15.768

$ws.Cells.Item(21, 1).Value = 2.590189054726355
$ws.Cells.Item(21, 3).Value = 10.636

$ws.Cells.Item(22, 1).Value = 7.075064676616908
$ws.Cells.Item(22, 2).Value = 0.188
$ws.Cells.Item(22, 3).Value = 28.58799999999999

$ws.Cells.Item(23, 1).Value = 2.183442786069638
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 9.875999999999994

$ws.Cells.Item(24, 1).Value = 3.559462686567157
$ws.Cells.Item(24, 3).Value = 14.508

$ws.Cells.Item(25, 1).Value = 2.829691542288548
$ws.Cells.Item(25, 3).Value = 11.50399999999999

$ws.Cells.Item(26, 1).Value = 2.893472636815908
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 12.028

$ws.Cells.Item(27, 1).Value = 6.98390049751243
$ws.Cells.Item(27, 2).Value = 0.1440000000000001
$ws.Cells.Item(27, 3).Value = 29.412

$ws.Cells.Item(28, 1).Value = 5.504079601990039
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 22.824

$ws.Cells.Item(29, 1).Value = 2.423144278606958
$ws.Cells.Item(29, 3).Value = 10.9

$ws.Cells.Item(30, 1).Value = 5.05980099502486
$ws.Cells.Item(30, 2).Value = 0.004
$ws.Cells.Item(30, 3).Value = 20.584

$ws.Cells.Item(31, 1).Value = 6.339223880597009
$ws.Cells.Item(31, 2).Value = 0.07200000000000002
$ws.Cells.Item(31, 3).Value = 25.22

$ws.Cells.Item(32, 1).Value = 6.307144278606956
$ws.Cells.Item(32, 2).Value = 0.036
$ws.Cells.Item(32, 3).Value = 25.772

$ws.Cells.Item(33, 1).Value = 6.571999999999987
$ws.Cells.Item(33, 2).Value = 0.024
$ws.Cells.Item(33, 3).Value = 25.87200000000001

$ws.Cells.Item(34, 1).Value = 5.89122388059701
$ws.Cells.Item(34, 2).Value = 0.044
$ws.Cells.Item(34, 3).Value = 24.248

$ws.Cells.Item(35, 1).Value = 3.584477611940291
$ws.Cells.Item(35, 3).Value = 14.692

$ws.Cells.Item(36, 1).Value = 5.008577114427847
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 20.532

$ws.Cells.Item(37, 1).Value = 2.782248756218895
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 11.5

$ws.Cells.Item(38, 1).Value = 6.991343283582082
$ws.Cells.Item(38, 2).Value = 0.1560000000000001
$ws.Cells.Item(38, 3).Value = 29.816

$ws.Cells.Item(39, 1).Value = 5.128398009950238
$ws.Cells.Item(39, 2).Value = 0.108
$ws.Cells.Item(39, 3).Value = 21.46

$ws.Cells.Item(40, 1).Value = 2.29603980099502
$ws.Cells.Item(40, 3).Value = 9.080000000000005

$ws.Cells.Item(41, 1).Value = 5.032915422885559
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 19.056

$ws.Cells.Item(42, 1).Value = 6.702069651741284
$ws.Cells.Item(42, 2).Value = 0.052
$ws.Cells.Item(42, 3).Value = 28.26400000000001

$ws.Cells.Item(43, 1).Value = 5.077213930348243
$ws.Cells.Item(43, 2).Value = 0.008
$ws.Cells.Item(43, 3).Value = 20.4

$ws.Cells.Item(44, 1).Value = 5.68553233830845
$ws.Cells.Item(44, 2).Value = 0.02
$ws.Cells.Item(44, 3).Value = 22.8

$ws.Cells.Item(45, 1).Value = 6.749054726368152
$ws.Cells.Item(45, 2).Value = 0.08000000000000002
$ws.Cells.Item(45, 3).Value = 26.50800000000001

$ws.Cells.Item(46, 1).Value = 6.234726368159194
$ws.Cells.Item(46, 2).Value = 0.004
$ws.Cells.Item(46, 3).Value = 24.19600000000001

$ws.Cells.Item(47, 1).Value = 7.003820895522382
$ws.Cells.Item(47, 2).Value = 0.028
$ws.Cells.Item(47, 3).Value = 26.90000000000001

$ws.Cells.Item(48, 1).Value = 5.147641791044768
$ws.Cells.Item(48, 2).Value = 0.03999999999999999
$ws.Cells.Item(48, 3).Value = 21.744

$ws.Cells.Item(49, 1).Value = 2.068935323383074
$ws.Cells.Item(49, 3).Value = 9.007999999999996

$ws.Cells.Item(50, 1).Value = 2.679263681592027
$ws.Cells.Item(50, 3).Value = 11.48

$ws.Cells.Item(51, 1).Value = 1.474169154228855
$ws.Cells.Item(51, 3).Value = 6.668000000000005

$ws.Cells.Item(52, 1).Value = 6.385014925373119
$ws.Cells.Item(52, 2).Value = 0.08000000000000002
$ws.Cells.Item(52, 3).Value = 23.608

$ws.Cells.Item(53, 1).Value = 5.297194029850736
$ws.Cells.Item(53, 2).Value = 0.08000000000000003
$ws.Cells.Item(53, 3).Value = 21.41999999999999

$ws.Cells.Item(54, 1).Value = 2.098646766169148
$ws.Cells.Item(54, 3).Value = 8.499999999999993

$ws.Cells.Item(55, 1).Value = 2.739721393034817
$ws.Cells.Item(55, 3).Value = 11.04

$ws.Cells.Item(56, 1).Value = 6.628457711442779
$ws.Cells.Item(56, 2).Value = 0.06000000000000001
$ws.Cells.Item(56, 3).Value = 27.48400000000001

$ws.Cells.Item(57, 1).Value = 6.537174129353231
$ws.Cells.Item(57, 2).Value = 0.02
$ws.Cells.Item(57, 3).Value = 25.732

$ws.Cells.Item(58, 1).Value = 2.933950248756209
$ws.Cells.Item(58, 2).Value = 0
$ws.Cells.Item(58, 3).Value = 12.516

$ws.Cells.Item(59, 1).Value = 6.17711442786069
$ws.Cells.Item(59, 2).Value = 0.004
$ws.Cells.Item(59, 3).Value = 26.26399999999999

$ws.Cells.Item(60, 1).Value = 2.108696517412924
$ws.Cells.Item(60, 3).Value = 9.468

$ws.Cells.Item(61, 1).Value = 7.004059701492524
$ws.Cells.Item(61, 2).Value = 0.108
$ws.Cells.Item(61, 3).Value = 29.236

$ws.Cells.Item(62, 1).Value = 6.271064676616906
$ws.Cells.Item(62, 2).Value = 0.008
$ws.Cells.Item(62, 3).Value = 24.108

$ws.Cells.Item(63, 1).Value = 6.875124378109445
$ws.Cells.Item(63, 2).Value = 0.08400000000000002
$ws.Cells.Item(63, 3).Value = 29.104

$ws.Cells.Item(64, 1).Value = 2.196378109452726
$ws.Cells.Item(64, 2).Value = 0
$ws.Cells.Item(64, 3).Value = 9.391999999999998

$ws.Cells.Item(65, 1).Value = 6.88346268656715
$ws.Cells.Item(65, 2).Value = 0.036
$ws.Cells.Item(65, 3).Value = 30.86000000000002

$ws.Cells.Item(66, 1).Value = 5.561910447761189
$ws.Cells.Item(66, 2).Value = 0
$ws.Cells.Item(66, 3).Value = 23.42

$ws.Cells.Item(67, 1).Value = 6.17522388059701
$ws.Cells.Item(67, 2).Value = 0.104
$ws.Cells.Item(67, 3).Value = 25.81999999999999

$ws.Cells.Item(68, 1).Value = 2.36461691542288
$ws.Cells.Item(68, 3).Value = 9.696

$ws.Cells.Item(69, 1).Value = 6.909930348258696
$ws.Cells.Item(69, 2).Value = 0.07600000000000003
$ws.Cells.Item(69, 3).Value = 27.36399999999999

$ws.Cells.Item(70, 1).Value = 2.242567164179098
$ws.Cells.Item(70, 3).Value = 9.736000000000001

$ws.Cells.Item(71, 1).Value = 3.928557213930342
$ws.Cells.Item(71, 3).Value = 15.46399999999998

$ws.Cells.Item(72, 1).Value = 4.355402985074621
$ws.Cells.Item(72, 3).Value = 17.84800000000001
